# Daily auto push: prepend two new daily-ranking rows (2026/01/16, 2026/01/17)
# ahead of the existing "2026/12/29" block, pushing the rest of the table
# down by two rows (638-679 -> 640-681).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the old row 638 ("2026/12/29" block),
# shifting everything from row 638 down to row 640 onward.
$ws.Rows.Item(638).Insert()
$ws.Rows.Item(639).Insert()

# Force column A/B on the new rows to behave as plain text (matching every
# other row in the sheet) instead of Excel's automatic date/number parsing.
$ws.Range("A638:B639").NumberFormat = "@"

$ws.Range("A638").Value = "2026/01/16"
$ws.Range("B638").Value = "金"
$ws.Range("C638").Value = 22
$ws.Range("D638").Value = 34

$ws.Range("A639").Value = "2026/01/17"
$ws.Range("B639").Value = "土"
$ws.Range("C639").Value = 2
$ws.Range("D639").Value = 34

# Drop the temporary "@" number format so the new cells carry no explicit
# style, same as the surrounding data rows.
$ws.Range("A638:D639").ClearFormats()
